$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for this product. It belongs
# chronologically right after the existing row 138 (old date 2022-05-24), so
# insert a fresh row there and push the rest of the table down by one.
$ws.Rows(138).Insert()

$ws.Range("A138").Value = 4
$ws.Range("B138").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C138").Value = "Los Lagos"
$ws.Range("D138").Value = 45041
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = 100112031
$ws.Range("G138").Value = "Poroto verde"
$ws.Range("H138").Value = "Magnum"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 40
$ws.Range("K138").Value = 35000
$ws.Range("L138").Value = 35000
$ws.Range("M138").Value = 35000
$ws.Range("N138").Value = "$/saco 25 kilos"
$ws.Range("O138").Value = "Región Metropolitana"
$ws.Range("P138").Value = 1400
$ws.Range("Q138").Value = 25
$ws.Range("R138").Value = "Hortaliza"
